$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (prices in column D, volume-change % in column E).
# Column D values are assigned with a leading apostrophe to force Excel to store
# them as literal text (quote-prefixed) instead of auto-converting number-like
# strings (e.g. "231.14", "16.00", "1.000") into numeric values.

$ws.Range("D2").Value = "'30.154.64"
$ws.Range("D3").Value = "'1.831.07"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'231.14"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4648"
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("D8").Value = "'0.2687"
$ws.Range("E8").Value = "  -6.88%  "
$ws.Range("D9").Value = "'0.06258"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("D10").Value = "'1.847.48"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "'16.00"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("D13").Value = "'4.894"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").Value = "'83.14"
$ws.Range("E14").Value = "  -5.53%  "
$ws.Range("D15").Value = "'0.6183"
$ws.Range("E15").Value = "  -7.47%  "
$ws.Range("D16").Value = "'30.081.10"
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'226.57"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").Value = "'0.000007272"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("D20").Value = "'12.37"
$ws.Range("E20").Value = "  -6.62%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'2.071.29"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'4.834"
$ws.Range("E23").Value = "  -8.50%  "
$ws.Range("D24").Value = "'5.860"
$ws.Range("E24").Value = "  -5.40%  "
$ws.Range("D25").Value = "'9.104"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("D26").Value = "'164.49"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").Value = "'17.61"
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "'1.840"
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").Value = "'0.1009"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").Value = "'4.041"
$ws.Range("E31").Value = "  -6.86%  "
$ws.Range("D32").Value = "'3.757"
$ws.Range("E32").Value = "  -6.73%  "
$ws.Range("D33").Value = "'0.04774"
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("D34").Value = "'1.123"
$ws.Range("E34").Value = "  -7.35%  "
$ws.Range("D35").Value = "'0.6981"
$ws.Range("E35").Value = "  -7.09%  "
$ws.Range("D36").Value = "'2.686"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "'0.01807"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").Value = "'2.606"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").Value = "'0.8930"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").Value = "'1.921"
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("D42").Value = "'102.73"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").Value = "'5.472"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").Value = "'0.3984"
$ws.Range("D45").Value = "'6.948"
$ws.Range("E45").Value = "  -6.45%  "
$ws.Range("D46").Value = "'0.1189"
$ws.Range("E46").Value = "  -7.12%  "
$ws.Range("D47").Value = "'59.49"
$ws.Range("E47").Value = "  -7.34%  "
$ws.Range("D48").Value = "'8.443"
$ws.Range("E48").Value = "  -6.29%  "
$ws.Range("D49").Value = "'0.05517"
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "'32.52"
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("D51").Value = "'1.357"
$ws.Range("E51").Value = "  -8.90%  "
